$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 / Row 12 new data ---------------------------------------------
# Writing in this specific order so that new shared-string entries are
# created with the same index ordering as the target workbook.

# 1) idx 47
$ws.Range("A11").Value = "102_AutomobileInsurance_006_SendQuote_001_MandatoryFields"
# 2) idx 48
$ws.Range("A12").Value = "102_AutomobileInsurance_006_SendQuote_002_FieldHintsAndErrors"
# 3) idx 49
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "102_AutomobileInsurance_006_SendQuote_002_EnterValuesInWrongFormat"
# 4) idx 50
$ws.Range("J12").NumberFormat = "@"
$ws.Range("J12").Value = "102_AutomobileInsurance_006_SendQuote_002_EnterValuesInWrongFormat Part 2"
# 5) idx 51
$ws.Range("F11").Value = "Choose Silver"
# 6) idx 52
$ws.Range("F12").Value = "Choose Gold"
# 7) idx 53
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = "Send Quote Page check error hint formatting"
# 8) idx 54
$ws.Range("K12").NumberFormat = "@"
$ws.Range("K12").Value = "Send Quote Page check error hint formatting Part 2"
# 9) idx 55
$ws.Range("G12").Value = "Send Quote Page check for hints regarding mandatory fields"
# 10) idx 56
$ws.Range("G11").Value = "Send Quote page check for open mandatory field"

# --- Reused cells (existing shared strings) --------------------------------
$ws.Range("B11").Value = "102_AutomobileInsurance_001_SmokeTest_FillPageVehicleData"
$ws.Range("C11").Value = "102_AutomobileInsurance_001_SmokeTest_FillPageInsurantData"
$ws.Range("D11").Value = "102_AutomobileInsurance_001_SmokeTest_FillPageProductData"
$ws.Range("E11").Value = "Goto price option page"

$ws.Range("B12").Value = "102_AutomobileInsurance_001_SmokeTest_FillPageVehicleData"
$ws.Range("C12").Value = "102_AutomobileInsurance_001_SmokeTest_FillPageInsurantData"
$ws.Range("D12").Value = "102_AutomobileInsurance_001_SmokeTest_FillPageProductData"
$ws.Range("E12").Value = "Goto price option page"

# --- Column width adjustments (columns A & B got a bit wider, J/K are new) -
$ws.Columns.Item(1).ColumnWidth = 58.944010416666664
$ws.Columns.Item(2).ColumnWidth = 56.830729166666664
$ws.Columns.Item(10).ColumnWidth = 69.49869791666667
$ws.Columns.Item(11).ColumnWidth = 39.276041666666664

# --- Extend header row to the two new columns (J1, K1) ---------------------
$ws.Range("J1").Value = "dlgAutomobileInsurance"
$ws.Range("K1").Value = "dlgAutomobileInsurance"

# --- Selection matches the authored workbook state -------------------------
[void]$ws.Range("G12").Select()
